$wb = $excel.ActiveWorkbook

# Add the new "Group B" sheet right after "Group A"
$groupA = $wb.Worksheets.Item("Group A")
$ws = $wb.Worksheets.Add($null, $groupA)
$ws.Name = "Group B"

# Header row (criteria / possible / score) - set font underline per-cell
# so we don't populate intermediate blank cells (B3, C3).
$ws.Range("A3").Value = "Criteria"
$ws.Range("A3").Font.Underline = 2

$ws.Range("D3").Value = "Possible"
$ws.Range("D3").Font.Underline = 2

$ws.Range("E3").Value = "Score"
$ws.Range("E3").Font.Underline = 2

# Rubric criteria rows
$ws.Range("A4").Value = "Are the correct developer’s name and date in a comment in the head element of the HTML page?"
$ws.Range("D4").Value = 5

$ws.Range("A5").Value = "Is a variable used for the favorite site name?"
$ws.Range("D5").Value = 5

$ws.Range("A6").Value = "Is a variable used for the site URL?"
$ws.Range("D6").Value = 10

$ws.Range("A7").Value = "Is the HTML that displays the link to the web site produced by a document.write statement?"
$ws.Range("D7").Value = 5

$ws.Range("A8").Value = "Is the caption produced by a document.write statement?"
$ws.Range("D8").Value = 5

$ws.Range("A9").Value = "Is a heading displayed for the page?"
$ws.Range("D9").Value = 5

$ws.Range("A10").Value = "Is the link displayed correctly?"
$ws.Range("D10").Value = 5

# Total row with formulas
$ws.Range("A12").Value = "Total"
$ws.Range("D12").Formula = "=SUM(D4:D10)"
$ws.Range("E12").Formula = "=SUM(E4:E10)"

# Title row (added last so the shared string for it lands at the end,
# matching the order new strings were appended in the source workbook)
$ws.Range("A1").Value = "Assignment Group B – Favorite web site"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A2").Font.Bold = $true

# Column C width to match Group A's sheet
$ws.Columns.Item(3).ColumnWidth = 16.83

# Match the zoom level used on the Group A sheet
$ws.Activate()
$excel.ActiveWindow.Zoom = 140
